$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the last used row based on worksheet dimension
$lastRow = $ws.Cells.SpecialCells(11).Row  # xlCellTypeLastCell = 11

# Column C holds the "Förändrad" (changed) date, stored as serial date 45202.
# All data rows (2..lastRow) need this date bumped by one day to 45203.
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45203
